# Rename the worksheet tab (drops the "-HW15.xpc" suffix, keeps the stem)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "AlphaFiberF"

# Update a few existing row-15 values to their higher-precision counterparts
$ws.Range("D15").Value = 0.9646606015518283
$ws.Range("J15").Value = 0.9646606015518283
$ws.Range("K15").Value = 0.9795612293667847
$ws.Range("L15").Value = 0.9967439034252417

# Append a new data row (row 16) for HKL index 14 / "HexGrid-60degTilt5degRes"
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C16").Value = 1.227435064709477
$ws.Range("D16").Value = 1.222548923474337
$ws.Range("E16").Value = 0.8672951964755753
$ws.Range("F16").Value = 1.227435064709477
$ws.Range("G16").Value = 1.00875899065177
$ws.Range("H16").Value = 0.936755875710869
$ws.Range("I16").Value = 0.9543713561677337
$ws.Range("J16").Value = 1.222548923474337
$ws.Range("K16").Value = 1.044922059974956
$ws.Range("L16").Value = 1.136178562342217
$ws.Range("M16").Value = 1.036194234531627

# Match style of column A data cells (A3:A15 use style index 1 / centered, bold, bordered)
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$excel.CutCopyMode = $false
